$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows/cells that no longer exist in the rerun (fewer anchor words this time)
$ws.Range("A9:H9").Clear()
$ws.Range("J31:Q32").Clear()

# Update cells whose value changed due to the larger dataset rerun
$ws.Range("A3").Value = 'crude'
$ws.Range("B3").Value = 0.8823529411764706
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 30
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = 'happy'
$ws.Range("K3").Value = 0.9615384615384616
$ws.Range("L3").Value = 25
$ws.Range("M3").Value = 25
$ws.Range("Q3").Value = 1
$ws.Range("B4").Value = 0.6111111111111112
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 22
$ws.Range("H4").Value = 14
$ws.Range("J4").Value = 'best'
$ws.Range("K4").Value = 0.9152542372881356
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 54
$ws.Range("Q4").Value = 5
$ws.Range("A5").Value = 'crisis'
$ws.Range("B5").Value = 0.5993150684931506
$ws.Range("C5").Value = 175
$ws.Range("D5").Value = 175
$ws.Range("H5").Value = 117
$ws.Range("J5").Value = 'love'
$ws.Range("K5").Value = 0.9130434782608695
$ws.Range("L5").Value = 42
$ws.Range("M5").Value = 42
$ws.Range("Q5").Value = 4
$ws.Range("A6").Value = 'emergency'
$ws.Range("B6").Value = 0.2533333333333334
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 19
$ws.Range("H6").Value = 56
$ws.Range("J6").Value = 'interesting'
$ws.Range("K6").Value = 0.9090909090909091
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = 30
$ws.Range("Q6").Value = 3
$ws.Range("B7").Value = 0.186046511627907
$ws.Range("C7").Value = 96
$ws.Range("D7").Value = 96
$ws.Range("H7").Value = 420
$ws.Range("J7").Value = 'great'
$ws.Range("K7").Value = 0.8571428571428571
$ws.Range("L7").Value = 96
$ws.Range("M7").Value = 96
$ws.Range("Q7").Value = 16
$ws.Range("B8").Value = 0.164021164021164
$ws.Range("C8").Value = 31
$ws.Range("D8").Value = 31
$ws.Range("H8").Value = 158
$ws.Range("J8").Value = 'nice'
$ws.Range("K9").Value = 0.7890625
$ws.Range("L9").Value = 101
$ws.Range("M9").Value = 101
$ws.Range("Q9").Value = 27
$ws.Range("J10").Value = 'thanks'
$ws.Range("K10").Value = 0.7804878048780488
$ws.Range("L10").Value = 64
$ws.Range("M10").Value = 64
$ws.Range("Q10").Value = 18
$ws.Range("J11").Value = 'positive'
$ws.Range("K11").Value = 0.7758620689655172
$ws.Range("L11").Value = 45
$ws.Range("M11").Value = 45
$ws.Range("Q11").Value = 13
$ws.Range("J12").Value = 'free'
$ws.Range("K12").Value = 0.725
$ws.Range("L12").Value = 87
$ws.Range("M12").Value = 87
$ws.Range("Q12").Value = 33
$ws.Range("J13").Value = 'confidence'
$ws.Range("K13").Value = 0.7222222222222222
$ws.Range("L13").Value = 26
$ws.Range("M13").Value = 26
$ws.Range("Q13").Value = 10
$ws.Range("J14").Value = 'special'
$ws.Range("K14").Value = 0.7222222222222222
$ws.Range("L14").Value = 26
$ws.Range("M14").Value = 26
$ws.Range("Q14").Value = 10
$ws.Range("J15").Value = 'good'
$ws.Range("K15").Value = 0.71875
$ws.Range("L15").Value = 115
$ws.Range("M15").Value = 115
$ws.Range("Q15").Value = 45
$ws.Range("J16").Value = 'support'
$ws.Range("K16").Value = 0.7075471698113207
$ws.Range("L16").Value = 75
$ws.Range("M16").Value = 75
$ws.Range("Q16").Value = 31
$ws.Range("J17").Value = 'safe'
$ws.Range("K17").Value = 0.6971830985915493
$ws.Range("L17").Value = 99
$ws.Range("M17").Value = 99
$ws.Range("Q17").Value = 43
$ws.Range("J18").Value = 'safety'
$ws.Range("K18").Value = 0.6470588235294118
$ws.Range("L18").Value = 33
$ws.Range("M18").Value = 33
$ws.Range("Q18").Value = 18
$ws.Range("J19").Value = 'well'
$ws.Range("K19").Value = 0.6063829787234043
$ws.Range("L19").Value = 57
$ws.Range("M19").Value = 57
$ws.Range("Q19").Value = 37
$ws.Range("J20").Value = 'relief'
$ws.Range("K20").Value = 0.58
$ws.Range("L20").Value = 29
$ws.Range("M20").Value = 29
$ws.Range("Q20").Value = 21
$ws.Range("K21").Value = 0.5714285714285714
$ws.Range("L21").Value = 36
$ws.Range("M21").Value = 36
$ws.Range("Q21").Value = 27
$ws.Range("J22").Value = 'fresh'
$ws.Range("K22").Value = 0.5416666666666666
$ws.Range("L22").Value = 26
$ws.Range("M22").Value = 26
$ws.Range("Q22").Value = 22
$ws.Range("J23").Value = 'hand'
$ws.Range("K23").Value = 0.5143603133159269
$ws.Range("L23").Value = 197
$ws.Range("M23").Value = 197
$ws.Range("Q23").Value = 186
$ws.Range("K24").Value = 0.4470588235294118
$ws.Range("L24").Value = 152
$ws.Range("M24").Value = 152
$ws.Range("Q24").Value = 188
$ws.Range("J25").Value = 'help'
$ws.Range("K25").Value = 0.4101694915254237
$ws.Range("L25").Value = 121
$ws.Range("M25").Value = 121
$ws.Range("Q25").Value = 174
$ws.Range("J26").Value = 'care'
$ws.Range("K26").Value = 0.4044943820224719
$ws.Range("L26").Value = 36
$ws.Range("M26").Value = 36
$ws.Range("Q26").Value = 53
$ws.Range("J27").Value = 'increase'
$ws.Range("K27").Value = 0.3461538461538461
$ws.Range("L27").Value = 27
$ws.Range("M27").Value = 27
$ws.Range("Q27").Value = 51
$ws.Range("J28").Value = 'please'
$ws.Range("K28").Value = 0.3305439330543933
$ws.Range("L28").Value = 79
$ws.Range("M28").Value = 79
$ws.Range("Q28").Value = 160
$ws.Range("J29").Value = 'protect'
$ws.Range("K29").Value = 0.3150684931506849
$ws.Range("L29").Value = 23
$ws.Range("M29").Value = 23
$ws.Range("Q29").Value = 50
$ws.Range("J30").Value = 'sure'
$ws.Range("K30").Value = 0.296875
$ws.Range("L30").Value = 19
$ws.Range("M30").Value = 19
$ws.Range("Q30").Value = 45
